# DSA Status workbook update - "Sheet Updated Till 30th Jan 2024"
# Adds 5 new LeetCode "Top Interview 150" rows (Merge Sorted Array, Remove
# Element, Remove Duplicates from Sorted Array, Majority Element,
# Linked List Cycle) below the existing row, each marked "Completed" with
# a green status fill, plus their hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New row data: Sl No, Problem Name, List Name, Platform, Link To Problem,
# Status, hyperlink-display-url
# ---------------------------------------------------------------------
$rows = @(
    @{ Row=3; SlNo=2; Name="Merge Sorted Array"; Problem="https://leetcode.com/problems/merge-sorted-array/description/?envType=study-plan-v2&envId=top-interview-150"; Link="https://leetcode.com/problems/merge-sorted-array/" },
    @{ Row=4; SlNo=3; Name="Remove Element"; Problem="https://leetcode.com/problems/remove-element/description/?envType=study-plan-v2&envId=top-interview-150"; Link="https://leetcode.com/problems/remove-element/" },
    @{ Row=5; SlNo=4; Name="Remove Duplicates from Sorted Array"; Problem="https://leetcode.com/problems/remove-duplicates-from-sorted-array/description/?envType=study-plan-v2&envId=top-interview-150"; Link="https://leetcode.com/problems/remove-duplicates-from-sorted-array/" },
    @{ Row=6; SlNo=5; Name="Majority Element"; Problem="https://leetcode.com/problems/majority-element/description/?envType=study-plan-v2&envId=top-interview-150"; Link="https://leetcode.com/problems/majority-element/" },
    @{ Row=7; SlNo=6; Name="Linked List Cycle"; Problem="https://leetcode.com/problems/linked-list-cycle/description/?envType=study-plan-v2&envId=top-interview-150"; Link="https://leetcode.com/problems/linked-list-cycle/" }
)

$listName = "Leetcode Top Interview 150"
$platform = "Leetcode"
$status = "Completed"

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Pull formatting (incl. the Hyperlink cell style) down from row 2 first
    # so values inherit the correct existing style indices instead of Excel
    # fabricating brand-new ones.
    $ws.Range("A2").Copy($ws.Range("A$rowNum"))
    $ws.Range("B2").Copy($ws.Range("B$rowNum"))
    $ws.Range("C2").Copy($ws.Range("C$rowNum"))
    $ws.Range("D2").Copy($ws.Range("D$rowNum"))
    $ws.Range("E2").Copy($ws.Range("E$rowNum"))
    $ws.Range("G2").Copy($ws.Range("G$rowNum"))

    $ws.Range("A$rowNum").Value = $r.SlNo
    $ws.Range("B$rowNum").Value = $r.Name
    $ws.Range("C$rowNum").Value = $listName
    $ws.Range("D$rowNum").Value = $platform
    $ws.Range("E$rowNum").Value = $r.Problem
    $ws.Range("G$rowNum").Value = $status

    # Hyperlink on the problem-name cell (mirrors B2 -> rId1 on row 2)
    $ws.Hyperlinks.Add($ws.Range("B$rowNum"), $r.Link, "", "", $r.Name)

    # Re-stamp the Hyperlink cell style (Hyperlinks.Add nudges the font),
    # keeping B$rowNum on the same style index as B2/B3.. instead of a dupe.
    $ws.Range("B2").Copy($ws.Range("B$rowNum"))
    $ws.Range("B$rowNum").Value = $r.Name
}

# Status fill for the new "Completed" rows: Green, Accent 6, Lighter 80%
# (theme 9 / tint 0.8) - distinct from the "In Progress" gold fill on G2.
$ws.Range("G3:G7").Interior.ThemeColor = 10
$ws.Range("G3:G7").Interior.TintAndShade = 0.79998168889431442

$ws.Range("D26").Select()
